$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D and E hold text-formatted values (e.g. "1.003", "  +5.90%  ").
# Force text number format so Excel does not auto-convert numeric-looking strings.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.146.44"
$ws.Range("E2").Value = "  +5.90%  "
$ws.Range("D3").Value = "1.914.14"
$ws.Range("E3").Value = "  +2.46%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.44%  "
$ws.Range("D5").Value = "330.32"
$ws.Range("E5").Value = "  +5.03%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("D7").Value = "0.5205"
$ws.Range("E7").Value = "  +2.76%  "
$ws.Range("D8").Value = "0.4083"
$ws.Range("E8").Value = "  +4.49%  "
$ws.Range("D9").Value = "0.08493"
$ws.Range("E9").Value = "  +2.07%  "
$ws.Range("D10").Value = "42.99"
$ws.Range("D11").Value = "1.126"
$ws.Range("E11").Value = "  +2.03%  "
$ws.Range("D12").Value = "23.01"
$ws.Range("E12").Value = "  +13.31%  "
$ws.Range("D13").Value = "6.444"
$ws.Range("E13").Value = "  +4.38%  "
$ws.Range("D14").Value = "1.919.35"
$ws.Range("E14").Value = "  +3.21%  "
$ws.Range("D15").Value = "7.385"
$ws.Range("E15").Value = "  +2.14%  "
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  -0.44%  "
$ws.Range("D17").Value = "95.20"
$ws.Range("E17").Value = "  +4.30%  "
$ws.Range("D18").Value = "0.00001113"
$ws.Range("E18").Value = "  +1.47%  "
$ws.Range("D19").Value = "0.06696"
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("D20").Value = "18.44"
$ws.Range("E20").Value = "  +4.74%  "
$ws.Range("D21").Value = "0.9994"
$ws.Range("D22").Value = "6.011"
$ws.Range("E22").Value = "  +1.82%  "
$ws.Range("D23").Value = "30.182.12"
$ws.Range("E23").Value = "  +5.94%  "
$ws.Range("D24").Value = "11.35"
$ws.Range("E24").Value = "  +2.78%  "
$ws.Range("D25").Value = "2.219"
$ws.Range("E25").Value = "  +1.17%  "
$ws.Range("D26").Value = "2.116.99"
$ws.Range("E26").Value = "  +2.05%  "
$ws.Range("D27").Value = "161.41"
$ws.Range("E27").Value = "  +2.35%  "
$ws.Range("D28").Value = "21.17"
$ws.Range("E28").Value = "  +3.25%  "
$ws.Range("D29").Value = "2.408"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("E30").Value = "  +2.32%  "
$ws.Range("D31").Value = "1.096"
$ws.Range("E31").Value = "  +5.83%  "
$ws.Range("D32").Value = "0.1069"
$ws.Range("E32").Value = "  +3.38%  "
$ws.Range("D33").Value = "5.999"
$ws.Range("E33").Value = "  +4.03%  "
$ws.Range("D34").Value = "3.605"
$ws.Range("E34").Value = "  -0.38%  "
$ws.Range("D35").Value = "0.02491"
$ws.Range("E35").Value = "  +2.11%  "
$ws.Range("D36").Value = "0.06577"
$ws.Range("E36").Value = "  +0.52%  "
$ws.Range("D37").Value = "0.2213"
$ws.Range("E37").Value = "  +2.58%  "
$ws.Range("D38").Value = "1.229"
$ws.Range("E38").Value = "  +4.02%  "
$ws.Range("D39").Value = "5.161"
$ws.Range("E39").Value = "  +2.86%  "
$ws.Range("D40").Value = "11.88"
$ws.Range("E40").Value = "  +7.24%  "
$ws.Range("D41").Value = "8.795"
$ws.Range("E41").Value = "  -1.81%  "
$ws.Range("D42").Value = "0.6516"
$ws.Range("E42").Value = "  +2.65%  "
$ws.Range("D43").Value = "1.238"
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("D44").Value = "0.6152"
$ws.Range("E44").Value = "  +3.04%  "
$ws.Range("D45").Value = "13.22"
$ws.Range("E45").Value = "  +1.83%  "
$ws.Range("D46").Value = "3.744"
$ws.Range("E46").Value = "  +1.94%  "
$ws.Range("D47").Value = "2.085"
$ws.Range("E47").Value = "  +4.60%  "
$ws.Range("D48").Value = "1.240"
$ws.Range("E48").Value = "  +2.69%  "
$ws.Range("D49").Value = "123.99"
$ws.Range("E49").Value = "  +1.78%  "
$ws.Range("E50").Value = "  +1.62%  "
$ws.Range("D51").Value = "79.54"
$ws.Range("E51").Value = "  +4.82%  "
